{"js": "// The document contains a paragraph built out of three runs:\n//   [<id>]  [p014r_1]  [</id>]\n// where the first/third run carry the \"tag\" formatting (Courier New,\n// color 7f6000, sz 18) and the middle run carries plain black text\n// formatting (color 000000). The authored edit collapses these three\n// runs into a single run - \"<id>p014r_1</id>\" - using the tag run's\n// formatting, i.e. the middle run's distinct formatting is dropped and\n// the text is merged into one contiguous run.\n//\n// Word (and Office.js) automatically merges adjacent runs that end up\n// with identical formatting once the \"odd\" run in the middle is\n// removed, so the approach is: locate the paragraph, find the\n// \"p014r_1\" sub-range inside it, delete that sub-range (which leaves\n// \"<id></id>\" as a single merged run carrying the tag formatting) and\n// then re-insert the \"p014r_1\" text at the gap (inheriting that same\n// run's formatting).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the exact paragraph \"<id>p014r_1</id>\" (works even if the\n// paragraph index in the document ever shifts).\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text === \"<id>p014r_1</id>\") {\n    target = p;\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not locate the '<id>p014r_1</id>' paragraph\");\n}\n\n// Search for the inner \"p014r_1\" text, scoped to this paragraph only,\n// so we don't touch the other (fig_p014r_1) occurrence elsewhere.\nconst paragraphRange = target.getRange();\nconst hits = paragraphRange.search(\"p014r_1\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length !== 1) {\n  throw new Error(\"Expected exactly one 'p014r_1' match inside the paragraph, found \" + hits.items.length);\n}\n\nconst hit = hits.items[0];\n// Remember the insertion point (collapsed caret) at the start of the hit\n// before deleting it, then delete the middle run's text (this merges the\n// surrounding \"<id>\" / \"</id>\" runs into a single run) and re-insert the\n// plain text there - it inherits the now-merged run's formatting.\nconst insertionPoint = hit.getRange(\"Start\");\nhit.delete();\ninsertionPoint.insertText(\"p014r_1\", \"Before\");\nawait context.sync();\n", "ps1": "# The document contains a paragraph built out of three runs:\n#   [<id>]  [p014r_1]  [</id>]\n# where the first/third run carry the \"tag\" formatting (Courier New,\n# color 7f6000, sz 18) and the middle run carries plain black text\n# formatting (color 000000). The authored edit collapses these three\n# runs into a single run - \"<id>p014r_1</id>\" - using the tag run's\n# formatting, i.e. the middle run's distinct formatting is dropped and\n# the text is merged into one contiguous run.\n#\n# Word automatically merges adjacent runs that end up with identical\n# formatting once the \"odd\" run in the middle is removed, so the\n# approach is: locate the paragraph, find the \"p014r_1\" sub-range\n# inside it, delete that sub-range (which leaves \"<id></id>\" as a\n# single merged run carrying the tag formatting) and then re-insert the\n# \"p014r_1\" text at the gap (inheriting that same run's formatting).\n\n$d = $word.ActiveDocument\n\n# Locate the target paragraph \"<id>p014r_1</id>\" robustly (not by a\n# hardcoded paragraph index) - trim the trailing paragraph mark (and any\n# cell-end mark) before comparing.\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"<id>p014r_1</id>\") {\n        $target = $p\n        break\n    }\n}\nif ($null -eq $target) {\n    throw \"Could not locate the '<id>p014r_1</id>' paragraph\"\n}\n\n# Scope Find to this paragraph's Range only, so the other unrelated\n# \"fig_p014r_1\" occurrence elsewhere in the document is left untouched.\n$pRange = $target.Range\n$found = $pRange.Find.Execute(\"p014r_1\")\nif (-not $found) {\n    throw \"Could not find 'p014r_1' inside the target paragraph\"\n}\n\n# $pRange now spans exactly the \"p014r_1\" run in the middle of the\n# paragraph (Find.Execute collapses/retargets the Range to the hit).\n$pStart = $pRange.Start\n$insertPoint = $d.Range($pStart, $pStart)\n\n# Deleting the middle run merges the surrounding \"<id>\" / \"</id>\" runs\n# into a single run (they share identical formatting); re-inserting the\n# plain text at the gap lands it inside that merged run.\n$pRange.Delete()\n$insertPoint.InsertAfter(\"p014r_1\")\n"}
